# Adding Board Outline mech layer to PCB/Outjob Template
# Insert a new column before the existing "Quantity" column (old column I) and
# populate its header with the new "Column=APSS Standard Part" field, shifting
# the old Quantity column from I to J.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at position I (9); this shifts the former column I
# (Quantity, and its styles/widths) one place to the right, to column J, while
# the new column I inherits the same per-row styling as the rest of the table.
$ws.Columns.Item(9).Insert()

# Set the header text for the newly inserted column I.
$ws.Range("I10").Value = "Column=APSS Standard Part"

# Set the new column's width (as close as this environment's pixel-rounded
# column-width model allows to the authored 22.453125 character-width units).
$ws.Columns.Item(9).ColumnWidth = 21.66

# Update the selection/active cell to match the saved view state, and scroll
# the sheet so the top-left visible cell reverts to the default.
$ws.Range("B14").Select()
